# Fill in the "floor type / cost per sqft" lookup table (G1:H4) and use it
# to compute per-run room costs (column D) for a set of sample runs, plus
# the final averaged overall cost in D8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lookup table (Table2, G1:H4). Populate in the order the new floor types
# first appear (hardwood, carpet, tile) so new shared strings line up.
$ws.Range("G2").Value = "hardwood"
$ws.Range("H2").Value = 1.39
$ws.Range("G4").Value = "carpet"
$ws.Range("H4").Value = 4.99
$ws.Range("G3").Value = "tile"
$ws.Range("H3").Value = 3.99

# Sample run 1: hardwood room, 20 x 10
$ws.Range("A3").Value = "hardwood"
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 10
$ws.Range("D3").Formula = "=(B3*C3)*H2"

# Sample run 2: carpet room, 20 x 10
$ws.Range("A4").Value = "carpet"
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = 10
$ws.Range("D4").Formula = "=(B4*C4)*H4"

# Sample run 3: hardwood room, 10 x 20
$ws.Range("A5").Value = "hardwood"
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 20
$ws.Range("D5").Formula = "=(B5*C5)*H2"

# Sample run 4: tile room, 20 x 15
$ws.Range("A6").Value = "tile"
$ws.Range("B6").Value = 20
$ws.Range("C6").Value = 15
$ws.Range("D6").Formula = "=(B6*C6)*H3"

# Sample run 5: carpet room, 15 x 10
$ws.Range("A7").Value = "carpet"
$ws.Range("B7").Value = 15
$ws.Range("C7").Value = 10
$ws.Range("D7").Formula = "=(B7*C7)*H4"

# Final overall (average) cost across the five sample runs
$ws.Range("D8").Formula = "=(D3+D4+D5+D6+D7)/5"

# Recalculate so cached formula results are written out
$excel.Calculate()

# Match the author's last selected cell
$ws.Range("D13").Select() | Out-Null
